$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "69298"
$ws.Range("A2").Style = "Normal"

# Row 3
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "07868"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = "Not being milked due to clinical mastitis"

# Row 4
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "44830"
$ws.Range("A4").Style = "Normal"

# Row 5
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "94638"
$ws.Range("A5").Style = "Normal"

# Row 6
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "47452"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = 2

# Row 7
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "62420"
$ws.Range("A7").Style = "Normal"

# Row 8
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "86112"
$ws.Range("A8").Style = "Normal"

# Row 9
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "23853"
$ws.Range("A9").Style = "Normal"

# Row 10
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "35117"
$ws.Range("A10").Style = "Normal"
$ws.Range("B10").Value = 1

# Row 11
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "25384"
$ws.Range("A11").Style = "Normal"

# Row 12
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "20750"
$ws.Range("A12").Style = "Normal"
$ws.Range("B12").Value = 1

# Row 13
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "88011"
$ws.Range("A13").Style = "Normal"
$ws.Range("B13").Value = 0

# Row 14
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "16027"
$ws.Range("A14").Style = "Normal"

# Row 15
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "67405"
$ws.Range("A15").Style = "Normal"

# Row 16
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "71725"
$ws.Range("A16").Style = "Normal"
$ws.Range("B16").Value = 1

# Row 17
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "28563"
$ws.Range("A17").Style = "Normal"
$ws.Range("B17").Value = 0

# Row 18
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "97054"
$ws.Range("A18").Style = "Normal"

# Row 19
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "07094"
$ws.Range("A19").Style = "Normal"
$ws.Range("B19").Value = 0

# Row 20
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "32717"
$ws.Range("A20").Style = "Normal"

# Row 21
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "61135"
$ws.Range("A21").Style = "Normal"

# Row 22
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "72616"
$ws.Range("A22").Style = "Normal"

# Row 23
$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = "35738"
$ws.Range("A23").Style = "Normal"

# Row 24
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = "92878"
$ws.Range("A24").Style = "Normal"

# Row 25
$ws.Range("A25").NumberFormat = "@"
$ws.Range("A25").Value = "75724"
$ws.Range("A25").Style = "Normal"
$ws.Range("B25").ClearContents()
$ws.Range("C25").Value = "Not being milked due to clinical mastitis"

# Row 26
$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = "93075"
$ws.Range("A26").Style = "Normal"
$ws.Range("B26").Value = 0

# Row 27
$ws.Range("A27").NumberFormat = "@"
$ws.Range("A27").Value = "81310"
$ws.Range("A27").Style = "Normal"
$ws.Range("B27").Value = 0

# Row 28
$ws.Range("A28").NumberFormat = "@"
$ws.Range("A28").Value = "57383"
$ws.Range("A28").Style = "Normal"
$ws.Range("B28").Value = 1

# Row 29
$ws.Range("A29").NumberFormat = "@"
$ws.Range("A29").Value = "55676"
$ws.Range("A29").Style = "Normal"

# Row 30
$ws.Range("A30").NumberFormat = "@"
$ws.Range("A30").Value = "87855"
$ws.Range("A30").Style = "Normal"

# Row 31
$ws.Range("A31").NumberFormat = "@"
$ws.Range("A31").Value = "70780"
$ws.Range("A31").Style = "Normal"

# Row 32
$ws.Range("A32").NumberFormat = "@"
$ws.Range("A32").Value = "48152"
$ws.Range("A32").Style = "Normal"

# Row 33
$ws.Range("A33").NumberFormat = "@"
$ws.Range("A33").Value = "80369"
$ws.Range("A33").Style = "Normal"
$ws.Range("B33").Value = 1

# Row 34
$ws.Range("A34").NumberFormat = "@"
$ws.Range("A34").Value = "13517"
$ws.Range("A34").Style = "Normal"

# Row 35
$ws.Range("A35").NumberFormat = "@"
$ws.Range("A35").Value = "88569"
$ws.Range("A35").Style = "Normal"
$ws.Range("B35").Value = 1

# Row 36
$ws.Range("A36").NumberFormat = "@"
$ws.Range("A36").Value = "72027"
$ws.Range("A36").Style = "Normal"
$ws.Range("B36").Value = 1

# Row 37
$ws.Range("A37").NumberFormat = "@"
$ws.Range("A37").Value = "29533"
$ws.Range("A37").Style = "Normal"

# Row 38
$ws.Range("A38").NumberFormat = "@"
$ws.Range("A38").Value = "85609"
$ws.Range("A38").Style = "Normal"

# Row 39
$ws.Range("A39").NumberFormat = "@"
$ws.Range("A39").Value = "97473"
$ws.Range("A39").Style = "Normal"

# Row 40
$ws.Range("A40").NumberFormat = "@"
$ws.Range("A40").Value = "30682"
$ws.Range("A40").Style = "Normal"

# Row 41
$ws.Range("A41").NumberFormat = "@"
$ws.Range("A41").Value = "33491"
$ws.Range("A41").Style = "Normal"

# Row 42
$ws.Range("A42").NumberFormat = "@"
$ws.Range("A42").Value = "58735"
$ws.Range("A42").Style = "Normal"
$ws.Range("B42").Value = 2

# Row 43
$ws.Range("A43").NumberFormat = "@"
$ws.Range("A43").Value = "67473"
$ws.Range("A43").Style = "Normal"
$ws.Range("B43").ClearContents()
$ws.Range("C43").Value = "Not being milked due to clinical mastitis"

# Row 44
$ws.Range("A44").NumberFormat = "@"
$ws.Range("A44").Value = "62919"
$ws.Range("A44").Style = "Normal"

# Row 45
$ws.Range("A45").NumberFormat = "@"
$ws.Range("A45").Value = "10979"
$ws.Range("A45").Style = "Normal"
$ws.Range("B45").Value = 0

# Row 46
$ws.Range("A46").NumberFormat = "@"
$ws.Range("A46").Value = "82048"
$ws.Range("A46").Style = "Normal"
$ws.Range("B46").Value = 0

# Row 47
$ws.Range("A47").NumberFormat = "@"
$ws.Range("A47").Value = "90999"
$ws.Range("A47").Style = "Normal"

# Row 48
$ws.Range("A48").NumberFormat = "@"
$ws.Range("A48").Value = "21049"
$ws.Range("A48").Style = "Normal"

# Row 49
$ws.Range("A49").NumberFormat = "@"
$ws.Range("A49").Value = "81953"
$ws.Range("A49").Style = "Normal"

# Row 50
$ws.Range("A50").NumberFormat = "@"
$ws.Range("A50").Value = "49845"
$ws.Range("A50").Style = "Normal"

# Row 51
$ws.Range("A51").NumberFormat = "@"
$ws.Range("A51").Value = "33975"
$ws.Range("A51").Style = "Normal"
